$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header labels: *_old -> *_FV2404, *_new -> *_FV2410
$suffixMap = @(
    "Segmentname", "Segmentgruppe", "Segment", "Datenelement", "Segment ID",
    "Code", "Qualifier", "Beschreibung", "Bedingungsausdruck", "Bedingung"
)

for ($i = 0; $i -lt $suffixMap.Length; $i++) {
    $col = $i + 1
    $ws.Cells.Item(1, $col).Value = "$($suffixMap[$i])_FV2404"
}

for ($i = 0; $i -lt $suffixMap.Length; $i++) {
    $col = $i + 12
    $ws.Cells.Item(1, $col).Value = "$($suffixMap[$i])_FV2410"
}

# Freeze the header row
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

# Turn the used range into an Excel Table
$range = $ws.Range("A1:U85")
$table = $ws.ListObjects.Add(1, $range, $null, 1)
$table.Name = "Table1"
